$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants_evaluated")
$r = $ws.Range("C4")
$r | Get-Member -MemberType Property | Out-String | Write-Output
